$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Problem Statement Approved" marker for a few students
# (Guttapati, Ma, Schaap) whose approval status reverted.
$ws.Range("J5").ClearContents()
$ws.Range("J7").ClearContents()
$ws.Range("J10").ClearContents()

# Reflect the new active selection left behind after the edit.
$ws.Range("J14").Select()
